$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $r = $d.Content
    $found = $r.Find.Execute($old, $true, $true, $false, $false, $false, `
                              $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Text not found: $old"
    }
    $r.Text = $new
}

Replace-Text "known as the ants puzzle, which I'm" "noto come il rompicapo delle formiche, che"
Replace-Text "probably going to discuss in a different" "probabilmente discuterò in un altro"
Replace-Text "video. Let me just finish writing down" "video. Lascia che finisca di scrivere"
Replace-Text "the title and, well, I can even draw a" "il titolo e, beh, posso anche disegnare una"
Replace-Text "little ant right here. okay, let's get" "piccola formica proprio qui. Okay, iniziamo"
Replace-Text "started! As I said I'm going to discuss" "! Come ho detto, discuterò"
Replace-Text "two puzzles in the first puzzle there" "due rompicapi. Nel primo ci"
Replace-Text "are two ants on a very high stool: a sort" "sono due formiche su uno sgabello molto alto: una sorta"
Replace-Text "of Mountain, flat on the top with two" "di Montagna, piatta in cima con due"
Replace-Text "steep cliffs to both the sides. The flat" "ripide scogliere su ambi i lati. Il picco"
Replace-Text "peak is one meter wide the two ants move" "piatto è largo un metro. Le due formiche si muovono"
Replace-Text "with a velocity, let's call it V, which is" "a una velocità, chiamiamola V, che è"
Replace-Text "the same for both of them and that is" "la stessa per entrambe e che è"
Replace-Text "equal to one centimeter per second. You" "pari a un centimetro al secondo. Puoi"
Replace-Text "can decide the direction towards each" "decidere la direzione verso cui ogni"
Replace-Text "ant moves if it is right or left and" "formica si muove, che sia destra o sinistra ed"
Replace-Text "where exactly to place the two ants on the" "esattamente dove posizionare le due formiche in"
Replace-Text "top of the mountain. Your purpose is to" "cima alla montagna. Il tuo scopo è"
Replace-Text "make the time the last ant takes before" "capire il tempo che l'ultima formica impiega prima"
Replace-Text "falling the longest possible. Ants cannot" "di cadere il più a lungo possibile. Le formiche non possono"
Replace-Text "be still: they must move to the right or" "rimanere ferme: devono muoversi a destra o"
Replace-Text "to the left but they must move and after" "a sinistra, ma devono muoversi e dopo"
Replace-Text "meeting each other they turn around and" "essersi incontrate, si voltano e"
Replace-Text "keep moving with the same but opposite" "continuano a muoversi ugualmente, ma alla velocità"
Replace-Text "velocity" "opposta"
Replace-Text "[Music]" "[Musica]"
Replace-Text "so again what are the precise positions" "Quindi, ancora, quali sono le posizioni precise"
Replace-Text "where I should place the two ants in" "in cui dovrei posizionare le due formiche,"
Replace-Text "order to get the longest time before the" "per ottenere il maggior tempo prima che"
Replace-Text "last ant falls? The second puzzle is" "l'ultima formica cada? Il secondo rompicapo è"
Replace-Text "basically the same but now we have three" "essenzialmente lo stesso, ma ora abbiamo tre"
Replace-Text "ants instead of two." "formiche invece di due."
Replace-Text "As before the ants velocity is one" "Come prima la velocità delle formiche è di un"
Replace-Text "centimeter per second, every ant turns" "centimetro al secondo, ogni formica si volta"
Replace-Text "around after meeting another ant and" "dopo averne incontrata un'altra e"
Replace-Text "the peak is one meter wide. So, what are" "il picco è largo un metro. Quindi, quali sono"
Replace-Text "now the precise positions" "ora le posizioni precise"
Replace-Text "I should place the three ants in order" "in cui dovrei posizionare le tre formiche per"
Replace-Text "to get the longest time before the last" "ottenere il maggior tempo prima che l'ultima"
Replace-Text "ant falls down? I hope you enjoyed this" "formica cada? Spero che questo video"
Replace-Text "video do your best and good luck" "ti sia piaciuto e buona fortuna."
